$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.788.15"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.942.82"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.48"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.34"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.940.26"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.63"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.99"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.961.08"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.426.39"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.939.20"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "445.27"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.80"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.680"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.49"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.00"
$ws.Range("E26").Value = "  -3.00%  "
$ws.Range("E27").Value = "  -6.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.37"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0992"
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.18"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.974"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.68"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.07"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  -5.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "43.69"
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.298"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.81"
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.119"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.41"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "382.04"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0351"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.724.84"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.02"
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.25"
$ws.Range("E51").Value = "  -0.01%  "
